# Natmi following Dr Hou advice
# Rewrites the LR-pair results table (rows 2..17, cols A..T) on Sheet1 with
# the recomputed ligand/receptor-expression statistics (now aggregated over
# 3 samples instead of 1), and extends the table from 13 to 17 data rows
# (4 sending clusters x 4 target clusters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs","Fn1","Itga2b","ECs",3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,3,1,0.6763496666666667,2.029049,0.1221603374409683,0.1221603374409683,13.49708695467245,121.473782592052,0.001433837655189034,0.001433837655189034),
  @("ECs","Fn1","Itga2b","FAPs",3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,3,1,3.724503,11.173509,0.672709052289864,0.6727090522898641,74.32537240934801,668.9283516841321,0.007895816190143049,0.007895816190143051),
  @("ECs","Fn1","Itga2b","M2",3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,2,0.6666666666666666,0.31493,0.94479,0.05688175357561716,0.05688175357561717,6.284674635213333,56.56207171692,0.0006676405933252705,0.0006676405933252707),
  @("ECs","Fn1","Itga2b","sCs",3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,3,1,0.8207906666666666,2.462372,0.1482488566935505,0.1482488566935505,16.37952015882844,147.415681429456,0.001740047527035143,0.001740047527035143),
  @("FAPs","Fn1","Itga2b","ECs",3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,3,1,0.6763496666666667,2.029049,0.1221603374409683,0.1221603374409683,1107.416624482385,9966.74962034147,0.1176443229192868,0.1176443229192868),
  @("FAPs","Fn1","Itga2b","FAPs",3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,3,1,3.724503,11.173509,0.672709052289864,0.6727090522898641,6098.29019427503,54884.61174847528,0.647840392685222,0.6478403926852222),
  @("FAPs","Fn1","Itga2b","M2",3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,2,0.6666666666666666,0.31493,0.94479,0.05688175357561716,0.05688175357561717,515.6485391159666,4640.8368520437,0.0547789530222843,0.05477895302228432),
  @("FAPs","Fn1","Itga2b","sCs",3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,3,1,0.8207906666666666,2.462372,0.1482488566935505,0.1482488566935505,1343.916134336795,12095.24520903116,0.1427684036784769,0.142768403678477),
  @("M2","Fn1","Itga2b","ECs",3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,3,1,0.6763496666666667,2.029049,0.1221603374409683,0.1221603374409683,11.83667332249356,106.530059902442,0.001257446734911002,0.001257446734911002),
  @("M2","Fn1","Itga2b","FAPs",3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,3,1,3.724503,11.173509,0.672709052289864,0.6727090522898641,65.18185410945802,586.6366869851221,0.006924471715344822,0.006924471715344825),
  @("M2","Fn1","Itga2b","M2",3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,2,0.6666666666666666,0.31493,0.94479,0.05688175357561716,0.05688175357561717,5.511533032646668,49.60379729382001,0.0005855073488499123,0.0005855073488499125),
  @("M2","Fn1","Itga2b","sCs",3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,3,1,0.8207906666666666,2.462372,0.1482488566935505,0.1482488566935505,14.36450916781956,129.280582510376,0.00152598662306148,0.001525986623061481),
  @("sCs","Fn1","Itga2b","ECs",3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,3,1,0.6763496666666667,2.029049,0.1221603374409683,0.1221603374409683,17.17665955112556,154.58993596013,0.001824730131581484,0.001824730131581484),
  @("sCs","Fn1","Itga2b","FAPs",3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,3,1,3.724503,11.173509,0.672709052289864,0.6727090522898641,94.58793754337,851.29143789033,0.01004837169915408,0.01004837169915409),
  @("sCs","Fn1","Itga2b","M2",3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,2,0.6666666666666666,0.31493,0.94479,0.05688175357561716,0.05688175357561717,7.998001121366666,71.98201009229999,0.0008496526111576754,0.0008496526111576756),
  @("sCs","Fn1","Itga2b","sCs",3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,3,1,0.8207906666666666,2.462372,0.1482488566935505,0.1482488566935505,20.84490100151555,187.60410901364,0.002214418864976923,0.002214418864976924)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $rowVals = $data[$i]
  for ($j = 0; $j -lt $rowVals.Length; $j++) {
    $c = $j + 1
    $ws.Cells.Item($r, $c).Value = $rowVals[$j]
  }
}

Write-Host "Wrote $($data.Length) data rows (rows 2..$($data.Length + 1))."
